$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated TPM-derived values to sheet1 (Rspo3-Sdc4 LR-pair table)
# Values below reflect the recomputed NATMI output after the TPM data update.

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1043256666666667
$ws.Cells.Item(2, 8).Value = 0.312977
$ws.Cells.Item(2, 9).Value = 0.02547563162231953
$ws.Cells.Item(2, 10).Value = 0.02547563162231953
$ws.Cells.Item(2, 13).Value = 0.5373756666666667
$ws.Cells.Item(2, 14).Value = 1.612127
$ws.Cells.Item(2, 15).Value = 0.007472820128982582
$ws.Cells.Item(2, 16).Value = 0.007472820128982581
$ws.Cells.Item(2, 17).Value = 0.05606207467544445
$ws.Cells.Item(2, 18).Value = 0.5045586720790001
$ws.Cells.Item(2, 19).Value = 0.0001903748127858146
$ws.Cells.Item(2, 20).Value = 0.0001903748127858146

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1043256666666667
$ws.Cells.Item(3, 8).Value = 0.312977
$ws.Cells.Item(3, 9).Value = 0.02547563162231953
$ws.Cells.Item(3, 10).Value = 0.02547563162231953
$ws.Cells.Item(3, 15).Value = 0.1537223653287423
$ws.Cells.Item(3, 16).Value = 0.1537223653287423
$ws.Cells.Item(3, 17).Value = 1.153245304396111
$ws.Cells.Item(3, 18).Value = 10.379207739565
$ws.Cells.Item(3, 19).Value = 0.003916174351226664
$ws.Cells.Item(3, 20).Value = 0.003916174351226663

# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1043256666666667
$ws.Cells.Item(4, 8).Value = 0.312977
$ws.Cells.Item(4, 9).Value = 0.02547563162231953
$ws.Cells.Item(4, 10).Value = 0.02547563162231953
$ws.Cells.Item(4, 13).Value = 30.561198
$ws.Cells.Item(4, 14).Value = 91.683594
$ws.Cells.Item(4, 15).Value = 0.4249882340167162
$ws.Cells.Item(4, 16).Value = 0.4249882340167161
$ws.Cells.Item(4, 17).Value = 3.188317355482
$ws.Cells.Item(4, 18).Value = 28.694856199338
$ws.Cells.Item(4, 19).Value = 0.01082684369362999
$ws.Cells.Item(4, 20).Value = 0.01082684369362999

# Row 5
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.1043256666666667
$ws.Cells.Item(5, 8).Value = 0.312977
$ws.Cells.Item(5, 9).Value = 0.02547563162231953
$ws.Cells.Item(5, 10).Value = 0.02547563162231953
$ws.Cells.Item(5, 13).Value = 29.75783666666667
$ws.Cells.Item(5, 14).Value = 89.27351
$ws.Cells.Item(5, 15).Value = 0.4138165805255589
$ws.Cells.Item(5, 16).Value = 0.4138165805255589
$ws.Cells.Item(5, 17).Value = 3.104506148807777
$ws.Cells.Item(5, 18).Value = 27.94055533927
$ws.Cells.Item(5, 19).Value = 0.01054223876467707
$ws.Cells.Item(5, 20).Value = 0.01054223876467707

# Row 6
$ws.Cells.Item(6, 9).Value = 0.9745243683776804
$ws.Cells.Item(6, 10).Value = 0.9745243683776804
$ws.Cells.Item(6, 13).Value = 0.5373756666666667
$ws.Cells.Item(6, 14).Value = 1.612127
$ws.Cells.Item(6, 15).Value = 0.007472820128982582
$ws.Cells.Item(6, 16).Value = 0.007472820128982581
$ws.Cells.Item(6, 17).Value = 2.144553615901889
$ws.Cells.Item(6, 18).Value = 19.300982543117
$ws.Cells.Item(6, 19).Value = 0.007282445316196767
$ws.Cells.Item(6, 20).Value = 0.007282445316196766

# Row 7
$ws.Cells.Item(7, 9).Value = 0.9745243683776804
$ws.Cells.Item(7, 10).Value = 0.9745243683776804
$ws.Cells.Item(7, 15).Value = 0.1537223653287423
$ws.Cells.Item(7, 16).Value = 0.1537223653287423
$ws.Cells.Item(7, 17).Value = 44.11532041727722
$ws.Cells.Item(7, 18).Value = 397.037883755495
$ws.Cells.Item(7, 19).Value = 0.1498061909775157
$ws.Cells.Item(7, 20).Value = 0.1498061909775156

# Row 8
$ws.Cells.Item(8, 9).Value = 0.9745243683776804
$ws.Cells.Item(8, 10).Value = 0.9745243683776804
$ws.Cells.Item(8, 13).Value = 30.561198
$ws.Cells.Item(8, 14).Value = 91.683594
$ws.Cells.Item(8, 15).Value = 0.4249882340167162
$ws.Cells.Item(8, 16).Value = 0.4249882340167161
$ws.Cells.Item(8, 17).Value = 121.963333553486
$ws.Cells.Item(8, 18).Value = 1097.670001981374
$ws.Cells.Item(8, 19).Value = 0.4141613903230862
$ws.Cells.Item(8, 20).Value = 0.4141613903230861

# Row 9
$ws.Cells.Item(9, 9).Value = 0.9745243683776804
$ws.Cells.Item(9, 10).Value = 0.9745243683776804
$ws.Cells.Item(9, 13).Value = 29.75783666666667
$ws.Cells.Item(9, 14).Value = 89.27351
$ws.Cells.Item(9, 15).Value = 0.4138165805255589
$ws.Cells.Item(9, 16).Value = 0.4138165805255589
$ws.Cells.Item(9, 17).Value = 118.7572869102455
$ws.Cells.Item(9, 18).Value = 1068.81558219221
$ws.Cells.Item(9, 19).Value = 0.4032743417608818
$ws.Cells.Item(9, 20).Value = 0.4032743417608818
